$wb = $excel.ActiveWorkbook

# The "Owner" column (column D) in the request-message sheets (publish,
# revise, query) was incorrectly listing the sender as "Publisher" -- it
# should read "Provider" instead. The response-message sheets
# (on_publish, on_revise, on_query) already had the correct values and
# are left untouched.
$sheetRows = @{
    "publish" = 17
    "revise"  = 16
    "query"   = 14
}

foreach ($sheetName in $sheetRows.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lastRow = $sheetRows[$sheetName]
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 4)
        if ($cell.Value2 -eq "Publisher") {
            $cell.Value = "Provider"
        }
    }
}
